$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.213.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.358.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.99%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.350.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.00%  "

$ws.Range("E10").Value = "  -4.40%  "

$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.335"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.781.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000163"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.323.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.368.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.47%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "63.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.74%  "

$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "548.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.67%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.470.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0906"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.05%  "

$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.46%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "150.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.363"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.70%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "138.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.579"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0495"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.00%  "

